$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.584.53"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.966.25"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.50"
$ws.Range("E5").Value = "  +6.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.02"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("E7").Value = "  -3.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.55"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0842"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.53"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.422.07"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.44"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.946.95"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.948"
$ws.Range("E17").Value = "  -7.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.616.79"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("E19").Value = "  -5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.89"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.80"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.172"
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "26.06"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("E29").Value = "  +7.44%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  -6.84%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.105"
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.96"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.06"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.12"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.31"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0432"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  -5.77%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.30"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -6.84%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.25"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.71"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.11"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.045.61"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.267"
$ws.Range("E47").Value = "  +7.78%  "
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.22"
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.244.84"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0329"
$ws.Range("E51").Value = "  -1.12%  "
